$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Teacher's name (appended after the existing "নাম:" label)
$ws.Range("A3").Value = "নাম: Dr. Sk. Imran Hossain"

# Designation (appended after the existing "পদবী: " label)
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"

# Year ("বর্ষ :") and Term ("টার্ম :") values
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"

# Department/Branch ("বিভাগ/শাখা:") value
$ws.Range("B5").Value = "সিএসই"

# Department ("বিভাগ :") appended with the value
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# Amount in words ("কথায়:") appended with the value
$ws.Range("A32").Value = "কথায়:তের হাজার ছয়শত পঁচান্ন টাকা মাত্র।"

# Widen column A so the name/designation text fits
$ws.Columns.Item(1).ColumnWidth = 14.33203125

# Leave the selection on the grand-total cell, matching where the user finished editing
$ws.Range("I32").Select()
